$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells that we are about to edit to be treated as plain text,
# matching the inlineStr/shared-string storage used by the source data (prices
# like "1.003" must not be reinterpreted as numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.665.13"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "1.892.86"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "311.65"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").Value = "0.3791"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").Value = "0.07312"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").Value = "0.9072"
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("D11").Value = "20.53"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("D12").Value = "0.07646"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "1.897.23"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "5.463"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").Value = "6.619"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "90.87"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "0.000008731"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "27.646.03"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "14.43"
$ws.Range("E21").Value = "  -4.35%  "
$ws.Range("D22").Value = "5.115"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "2.128.39"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").Value = "10.74"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").Value = "153.95"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").Value = "1.857"
$ws.Range("E26").Value = "  -5.40%  "
$ws.Range("D27").Value = "18.36"
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("D28").Value = "2.153"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "115.21"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "4.859"
$ws.Range("E30").Value = "  -3.73%  "
$ws.Range("D31").Value = "0.08940"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("D33").Value = "1.222"
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("D34").Value = "0.7611"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").Value = "4.619"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").Value = "2.548"
$ws.Range("E36").Value = "  -8.86%  "
$ws.Range("D37").Value = "0.02031"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").Value = "0.05285"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("D40").Value = "2.989"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "0.5445"
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("D42").Value = "6.866"
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("D43").Value = "8.482"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").Value = "112.46"
$ws.Range("E44").Value = "  +4.89%  "
$ws.Range("D45").Value = "0.1516"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D47").Value = "0.4776"
$ws.Range("E47").Value = "  -4.04%  "
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "1.629"
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("D50").Value = "67.15"
$ws.Range("E50").Value = "  -3.60%  "
$ws.Range("E51").Value = "  -1.13%  "

# Restore the original (default) style on column D now that the literal text
# values are committed, so we do not leave a stray number-format override.
$ws.Range("D2:D51").Style = "Normal"
